$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.2
$ws.Range("H2").Value = 3.35
$ws.Range("I2").Value = 3.4
$ws.Range("N2").Value = 4.6
$ws.Range("P2").Value = 2.22
$ws.Range("Q2").Value = 1.77
$ws.Range("S2").Value = 2.96
$ws.Range("T2").Value = 1.67
$ws.Range("AC2").Value = 8.800000000000001
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 38
$ws.Range("AN2").Value = 14
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 13
$ws.Range("AW2").Value = 27
$ws.Range("AZ2").Value = 14.5
$ws.Range("BA2").Value = 24
$ws.Range("BC2").Value = 19
$ws.Range("BD2").Value = 27
$ws.Range("BG2").Value = 26
$ws.Range("BH2").Value = "2026-02-23 08:31:48"

# Row 3
$ws.Range("H3").Value = 2.62
$ws.Range("K3").Value = 3.6
$ws.Range("P3").Value = 1.92
$ws.Range("BH3").Value = "2026-02-23 08:31:48"

# Row 4
$ws.Range("F4").Value = 15.5
$ws.Range("G4").Value = 22
$ws.Range("H4").Value = 1.17
$ws.Range("I4").Value = 1.24
$ws.Range("J4").Value = 7.8
$ws.Range("K4").Value = 11
$ws.Range("P4").Value = 4.2
$ws.Range("Q4").Value = 1.22
$ws.Range("BH4").Value = "2026-02-23 08:31:48"

# Row 5
$ws.Range("P5").Value = 1.78
$ws.Range("Q5").Value = 2.08
$ws.Range("BH5").Value = "2026-02-23 08:31:48"

# Row 6
$ws.Range("F6").Value = 2.24
$ws.Range("G6").Value = 2.74
$ws.Range("H6").Value = 3.75
$ws.Range("K6").Value = 3.25
$ws.Range("P6").Value = 1.43
$ws.Range("Q6").Value = 2.84
$ws.Range("BH6").Value = "2026-02-23 08:31:48"

# Row 7
$ws.Range("P7").Value = 1.53
$ws.Range("Q7").Value = 2.52
$ws.Range("BH7").Value = "2026-02-23 08:31:48"

# Row 8
$ws.Range("G8").Value = 2.82
$ws.Range("I8").Value = 4.5
$ws.Range("K8").Value = 3.45
$ws.Range("Q8").Value = 2.22
$ws.Range("T8").Value = 1.92
$ws.Range("U8").Value = 1.84
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 13.5
$ws.Range("AA8").Value = 95
$ws.Range("AB8").Value = 9.800000000000001
$ws.Range("AC8").Value = 8.6
$ws.Range("AH8").Value = 970
$ws.Range("AJ8").Value = 44
$ws.Range("AM8").Value = 180
$ws.Range("AN8").Value = 38
$ws.Range("AP8").Value = 3.2
$ws.Range("AQ8").Value = 3.3
$ws.Range("AR8").Value = 3.85
$ws.Range("AS8").Value = 4.2
$ws.Range("AT8").Value = 7
$ws.Range("AU8").Value = 6
$ws.Range("AV8").Value = 3.6
$ws.Range("AW8").Value = 4.1
$ws.Range("AX8").Value = 3.5
$ws.Range("AY8").Value = 3.4
$ws.Range("AZ8").Value = 3.75
$ws.Range("BA8").Value = 4.2
$ws.Range("BB8").Value = 4
$ws.Range("BC8").Value = 3.95
$ws.Range("BD8").Value = 4.1
$ws.Range("BE8").Value = 4.3
$ws.Range("BF8").Value = 3.95
$ws.Range("BG8").Value = 4.2
$ws.Range("BH8").Value = "2026-02-23 08:31:48"

# Row 9
$ws.Range("G9").Value = 3.2
$ws.Range("I9").Value = 3.2
$ws.Range("BH9").Value = "2026-02-23 08:31:48"

# Row 10
$ws.Range("G10").Value = 2.44
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.95
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 4.4
$ws.Range("P10").Value = 1.81
$ws.Range("Q10").Value = 1.88
$ws.Range("BH10").Value = "2026-02-23 08:31:48"

# Row 11
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 1.95
$ws.Range("I11").Value = 2.26
$ws.Range("BH11").Value = "2026-02-23 08:31:48"

# Row 12
$ws.Range("F12").Value = 2.3
$ws.Range("H12").Value = 2.96
$ws.Range("I12").Value = 3.45
$ws.Range("BH12").Value = "2026-02-23 08:31:48"

# Row 13
$ws.Range("F13").Value = 2.46
$ws.Range("H13").Value = 2.84
$ws.Range("I13").Value = 3.25
$ws.Range("P13").Value = 1.86
$ws.Range("Q13").Value = 1.93
$ws.Range("BH13").Value = "2026-02-23 08:31:48"

# Row 14
$ws.Range("Q14").Value = 3.2
$ws.Range("BH14").Value = "2026-02-23 08:31:48"

# Row 15
$ws.Range("F15").Value = 2.78
$ws.Range("G15").Value = 4.2
$ws.Range("H15").Value = 2.54
$ws.Range("I15").Value = 3
$ws.Range("BH15").Value = "2026-02-23 08:31:48"

# Row 16
$ws.Range("F16").Value = 1.54
$ws.Range("G16").Value = 1.55
$ws.Range("H16").Value = 5.9
$ws.Range("J16").Value = 5.3
$ws.Range("K16").Value = 5.4
$ws.Range("N16").Value = 6.6
$ws.Range("P16").Value = 2.94
$ws.Range("Q16").Value = 1.49
$ws.Range("S16").Value = 2.24
$ws.Range("AI16").Value = 55
$ws.Range("AK16").Value = 13.5
$ws.Range("AL16").Value = 25
$ws.Range("AN16").Value = 5.2
$ws.Range("AR16").Value = 28
$ws.Range("BD16").Value = 22
$ws.Range("BF16").Value = 5
$ws.Range("BG16").Value = 42
$ws.Range("BH16").Value = "2026-02-23 08:31:48"

# Row 17
$ws.Range("G17").Value = 1.51
$ws.Range("H17").Value = 7
$ws.Range("I17").Value = 7.2
$ws.Range("K17").Value = 5.3
$ws.Range("P17").Value = 2.86
$ws.Range("Q17").Value = 1.51
$ws.Range("R17").Value = 1.76
$ws.Range("S17").Value = 2.26
$ws.Range("U17").Value = 2.34
$ws.Range("AH17").Value = 19.5
$ws.Range("AL17").Value = 26
$ws.Range("AN17").Value = 5
$ws.Range("AU17").Value = 11
$ws.Range("AW17").Value = 34
$ws.Range("AX17").Value = 10.5
$ws.Range("BA17").Value = 32
$ws.Range("BF17").Value = 4.7
$ws.Range("BG17").Value = 32
$ws.Range("BH17").Value = "2026-02-23 08:31:48"

# Row 18
$ws.Range("I18").Value = 12
$ws.Range("N18").Value = 9.199999999999999
$ws.Range("P18").Value = 3.7
$ws.Range("Q18").Value = 1.34
$ws.Range("R18").Value = 2.12
$ws.Range("S18").Value = 1.85
$ws.Range("T18").Value = 1.69
$ws.Range("X18").Value = 50
$ws.Range("AB18").Value = 16.5
$ws.Range("AC18").Value = 18.5
$ws.Range("AI18").Value = 95
$ws.Range("AJ18").Value = 12.5
$ws.Range("AM18").Value = 85
$ws.Range("AO18").Value = 100
$ws.Range("AP18").Value = 44
$ws.Range("AQ18").Value = 48
$ws.Range("AS18").Value = 100
$ws.Range("AW18").Value = 42
$ws.Range("AY18").Value = 11
$ws.Range("BA18").Value = 65
$ws.Range("BF18").Value = 3.1
$ws.Range("BG18").Value = 42
$ws.Range("BH18").Value = "2026-02-23 08:31:48"

# Row 19
$ws.Range("BH19").Value = "2026-02-23 08:31:48"

# Row 20
$ws.Range("BH20").Value = "2026-02-23 08:31:48"

# Row 21
$ws.Range("F21").Value = 1.87
$ws.Range("G21").Value = 2
$ws.Range("I21").Value = 7.6
$ws.Range("J21").Value = 3
$ws.Range("P21").Value = 1.41
$ws.Range("Q21").Value = 2.92
$ws.Range("BH21").Value = "2026-02-23 08:31:48"

# Row 22
$ws.Range("H22").Value = 2.48
$ws.Range("BH22").Value = "2026-02-23 08:31:48"

# Row 23
$ws.Range("BH23").Value = "2026-02-23 08:31:48"
